# Trained the XGBoost models.
#
# Fills in the Accuracy (col E) / Optimal Parameters (col F) results that came
# out of (re-)training the K-Nearest Neighbors and Boosted Trees (XGBoost)
# models in the "model_analysis" comparison table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K-Nearest Neighbors / Distance (3 classes), RSSI row ---
$ws.Range("E10").Value = 0.58615555555555499
$ws.Range("F10").Value = "n_neighbors=1381, metric='manhattan'"

# --- K-Nearest Neighbors / Distance (Binary) rows ---
$ws.Range("E14").Value = 0.67968333333333297
$ws.Range("F14").Value = "n_neighbors=1081, metric='manhattan'"

$ws.Range("E15").Value = 0.97416666666666596
$ws.Range("F15").Value = "n_neighbors=19, metric='manhattan', n_components=2"

# --- Boosted Trees (XGBoost) / Distance (3 classes) rows ---
$ws.Range("E18").Value = 0.59760000000000002
$ws.Range("F18").Value = "n_estimators=3, max_depth=3, min_child_weight=1, gamma=0, colsample_bytree=0.01, subsample=0.1, reg_alpha=0, reg_lambda=0"

$ws.Range("F19").Value = "n_estimators=2139, max_depth=3, min_child_weight=3, gamma=0.1, colsample_bytree=0.01, subsample=0.2, reg_alpha=0.01, reg_lambda=0"

$ws.Range("F20").Value = "n_estimators=37, max_depth=5, min_child_weight=1, gamma=0.1, colsample_bytree=0.01, subsample=0.8, reg_alpha=5, reg_lambda=0"

$ws.Range("F21").Value = "n_estimators=396, max_depth=5, min_child_weight=3, gamma=0, colsample_bytree=0.7, subsample=0.8, reg_alpha=4, reg_lamba=0"

# --- Boosted Trees (XGBoost) / Distance (Binary) rows ---
$ws.Range("E22").Value = 0.6794
$ws.Range("F22").Value = "n_estimators=485, max_depth=3, min_child_weight=1, gamma=0, colsample_bytree=0.01, subsample=0.03, reg_alpha=0, reg_lambda=0"

$ws.Range("E23").Value = 0.97340000000000004
$ws.Range("F23").Value = "n_estimators=2815, max_depth=4, min_child_weight=2, gamma=0.2, colsample_bytree=0.01, subsample=0.76, reg_alpha=0.01, reg_lambda=0"

$ws.Range("E24").Value = 0.76989300000000005
$ws.Range("F24").Value = "n_estimators=733, max_depth=4, min_child_weight=3, gamma=0.1, colsample_bytree=0.01, subsample=0.83, reg_alpha=0.01, reg_lambda=0"

$ws.Range("E25").Value = 0.9819
$ws.Range("F25").Value = "n_estimators=2445, max_depth=5, min_child_weight=3, gamma=0.3, colsample_bytree=0.67, subsample=0.78, reg_alpha=0, reg_lambda=0"

# --- Minor workbook-state cosmetics: the author widened column F a bit and
# left the selection on F27 after scrolling the view over to column C. ---
$ws.Columns.Item(6).ColumnWidth = 125.8

$ws.Range("C1").Select() | Out-Null
$ws.Range("F27").Select() | Out-Null
